$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: fix number formats for F/G cells that are currently empty (style 1) so they
# adopt the existing OF_CDG (style 2) / OF_DATA (style 3) number formats used elsewhere. ---
$ws.Range("F350").Copy()
$ws.Range("F348").PasteSpecial(-4122)
$ws.Range("F351").PasteSpecial(-4122)
$ws.Range("F352").PasteSpecial(-4122)
$ws.Range("F353").PasteSpecial(-4122)
$ws.Range("F354").PasteSpecial(-4122)
$ws.Range("F355").PasteSpecial(-4122)
$ws.Range("G350").Copy()
$ws.Range("G348").PasteSpecial(-4122)
$ws.Range("G351").PasteSpecial(-4122)
$ws.Range("G352").PasteSpecial(-4122)
$ws.Range("G353").PasteSpecial(-4122)
$ws.Range("G354").PasteSpecial(-4122)
$ws.Range("G355").PasteSpecial(-4122)

# --- Step 2: fix O column cells that hold zero-padded numeric-looking codes; copy the
# exact text value (and keep destination formatting) from another row already containing
# the same FORNECEDOR_CDG text, so Excel does not reinterpret it as a number. ---
$ws.Range("O2").Copy()
$ws.Range("O347").PasteSpecial(-4163)
$ws.Range("O2").Copy()
$ws.Range("O348").PasteSpecial(-4163)
$ws.Range("O102").Copy()
$ws.Range("O349").PasteSpecial(-4163)
$ws.Range("O63").Copy()
$ws.Range("O350").PasteSpecial(-4163)
$ws.Range("O56").Copy()
$ws.Range("O351").PasteSpecial(-4163)
$ws.Range("O63").Copy()
$ws.Range("O352").PasteSpecial(-4163)
$ws.Range("O56").Copy()
$ws.Range("O353").PasteSpecial(-4163)
$ws.Range("O63").Copy()
$ws.Range("O354").PasteSpecial(-4163)
$ws.Range("O63").Copy()
$ws.Range("O355").PasteSpecial(-4163)
$excel.CutCopyMode = 0

# --- Step 3: set remaining simple values directly. ---
# Row 347
$ws.Range("F347").Value = 81195
$ws.Range("G347").Value = 46002
$ws.Range("H347").Value = "E.04.0892"
$ws.Range("I347").Value = "DISCO DE CORTE AÇO INOX 4 1/2''"
$ws.Range("J347").Value = "Apto"
$ws.Range("K347").Value = "UN"
$ws.Range("L347").Value = 20
$ws.Range("M347").Value = 3.5
$ws.Range("N347").Value = 70
$ws.Range("P347").Value = "DMC MATERIAIS"

# Row 348
$ws.Range("F348").Value = 81195
$ws.Range("G348").Value = 46002
$ws.Range("H348").Value = "E.04.0716"
$ws.Range("I348").Value = "DISCO DE CORTE DIAMANTADO PARA PORCELANATO  CONTINUO ( 4 1/2'' )"
$ws.Range("J348").Value = "Apto"
$ws.Range("K348").Value = "UN"
$ws.Range("L348").Value = 3
$ws.Range("M348").Value = 35
$ws.Range("N348").Value = 105
$ws.Range("P348").Value = "DMC MATERIAIS"

# Row 349
$ws.Range("F349").Value = 81209
$ws.Range("G349").Value = 46003
$ws.Range("H349").Value = "C.04.0012"
$ws.Range("I349").Value = "PANO MULTIUSO"
$ws.Range("J349").Value = "Apto"
$ws.Range("K349").Value = "UN"
$ws.Range("L349").Value = 2
$ws.Range("M349").Value = 19.9
$ws.Range("N349").Value = 39.8
$ws.Range("P349").Value = "VILE EMBALAGENS"

# Row 350
$ws.Range("F350").Value = 81197
$ws.Range("G350").Value = 46002
$ws.Range("H350").Value = "E.02.0040"
$ws.Range("I350").Value = "LONA PLÁSTICA  TRANSPARENTE  - 4 X 100 M"
$ws.Range("J350").Value = "Apto"
$ws.Range("K350").Value = "RL"
$ws.Range("L350").Value = 4
$ws.Range("M350").Value = 305
$ws.Range("N350").Value = 1220
$ws.Range("P350").Value = "CASA PEDROSO2648864-"

# Row 351
$ws.Range("F351").Value = 81206
$ws.Range("G351").Value = 46003
$ws.Range("H351").Value = "E.02.0012"
$ws.Range("I351").Value = "FITA CREPE LARANJA 45 MM X  50M"
$ws.Range("J351").Value = "Apto"
$ws.Range("K351").Value = "UN"
$ws.Range("L351").Value = 20
$ws.Range("M351").Value = 23.43
$ws.Range("N351").Value = 468.6
$ws.Range("P351").Value = "BELLOOBRAS"

# Row 352
$ws.Range("F352").Value = 81197
$ws.Range("G352").Value = 46002
$ws.Range("H352").Value = "E.02.0065"
$ws.Range("I352").Value = "PAPELÃO COM LOGO OSBORNE"
$ws.Range("J352").Value = "Apto"
$ws.Range("K352").Value = "UN"
$ws.Range("L352").Value = 10
$ws.Range("M352").Value = 90
$ws.Range("N352").Value = 900
$ws.Range("P352").Value = "CASA PEDROSO2648864-"

# Row 353
$ws.Range("F353").Value = 81206
$ws.Range("G353").Value = 46003
$ws.Range("H353").Value = "E.04.0405"
$ws.Range("I353").Value = "ESTILETE AUTOMÁTICO EMBORRACHADO 8 LÂMINAS - 18MM"
$ws.Range("J353").Value = "Apto"
$ws.Range("K353").Value = "UN"
$ws.Range("L353").Value = 6
$ws.Range("M353").Value = 10
$ws.Range("N353").Value = 60
$ws.Range("P353").Value = "BELLOOBRAS"

# Row 354
$ws.Range("F354").Value = 81197
$ws.Range("G354").Value = 46002
$ws.Range("H354").Value = "J.03.0015"
$ws.Range("I354").Value = "AREIA  - SACO GRANDE 20KG"
$ws.Range("J354").Value = "Apto"
$ws.Range("K354").Value = "SC"
$ws.Range("L354").Value = 20
$ws.Range("M354").Value = 5.8
$ws.Range("N354").Value = 116
$ws.Range("P354").Value = "CASA PEDROSO2648864-"

# Row 355
$ws.Range("F355").Value = 81197
$ws.Range("G355").Value = 46002
$ws.Range("H355").Value = "J.05.0002"
$ws.Range("I355").Value = "CIMENTO CP II - E-32 - 25 KG"
$ws.Range("J355").Value = "Apto"
$ws.Range("K355").Value = "SC"
$ws.Range("L355").Value = 5
$ws.Range("M355").Value = 37.9
$ws.Range("N355").Value = 189.5
$ws.Range("P355").Value = "CASA PEDROSO2648864-"
